# ============================================================================
# Applies the "additional scraping" edit:
#   1. Inserts a new "Player Info" sheet at the front with ID/NAME/
#      BATTING_HAND/BOWL_STYLE for the player.
#   2. Renames MATCH_CARD_LINK -> MATCH_CODE on the existing "ODI Batting"
#      and "ODI Bowling" sheets, replacing the full scorecard URL with just
#      the numeric match code, and drops the stray empty INNING_NUMBER
#      cells on "ODI Batting".
#   3. Appends a new "ODI Batting Extra" sheet with additional per-match
#      batting detail (batting position, boundary counts, % of total runs,
#      man-of-the-match flag).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Helper: stamp a 1-row header range with the same look as the existing
# bold / thin-bordered / center-top aligned header cells already used on the
# "ODI Batting" / "ODI Bowling" sheets.
# ----------------------------------------------------------------------------
function Format-HeaderRange($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous (thin box border)
}

# ----------------------------------------------------------------------------
# 1. Insert "Player Info" before "ODI Batting", append "ODI Batting Extra"
#    after "ODI Bowling", so the final tab order is:
#       Player Info | ODI Batting | ODI Bowling | ODI Batting Extra
# ----------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$extraSheet = $wb.Worksheets.Add($null, $bowlingSheet)
$extraSheet.Name = "ODI Batting Extra"

# ----------------------------------------------------------------------------
# 2. Populate "Player Info"
# ----------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Item("Player Info")

$piHeader = $playerInfo.Range("A1:D1")
Format-HeaderRange $piHeader
$playerInfo.Cells.Item(1,1).Value = "ID"
$playerInfo.Cells.Item(1,2).Value = "NAME"
$playerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"

# Force column A's data cell to stay text (so "4675" isn't coerced to a number)
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Cells.Item(2,1).Value = "4675"
$playerInfo.Cells.Item(2,2).Value = "Mohammed Siraj"
$playerInfo.Cells.Item(2,3).Value = "Right Handed"
$playerInfo.Cells.Item(2,4).Value = "Right Arm Medium Fast"

# ----------------------------------------------------------------------------
# 3. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (column D), drop the
#    leftover empty INNING_NUMBER cells (column B) for the "did not bat"
#    rows.
# ----------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Item(1,4).Value = "MATCH_CODE"

$battingCodes = @{
    2  = "4235"; 3  = "4533"; 4  = "4535"; 5  = "4536"; 6  = "4618";
    7  = "4621"; 8  = "4623"; 9  = "4624"; 10 = "4637"; 11 = "4640";
    12 = "4656"; 13 = "4657"; 14 = "4658"; 15 = "4679"; 16 = "4682";
    17 = "4685"; 18 = "4687"; 19 = "4689"; 20 = "4691"; 21 = "4692";
    22 = "4695"; 23 = "4725"; 24 = "4728"; 25 = "4732"
}
$battingSheet.Range("D2:D25").NumberFormat = "@"
foreach ($r in $battingCodes.Keys) {
    $battingSheet.Cells.Item($r, 4).Value = $battingCodes[$r]
}

$emptyInningRows = @(2, 3, 6, 9, 10, 11, 12, 13, 14, 19, 20, 21, 22, 23)
foreach ($r in $emptyInningRows) {
    $battingSheet.Cells.Item($r, 2).ClearContents()
}

# ----------------------------------------------------------------------------
# 4. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (column B)
# ----------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1,2).Value = "MATCH_CODE"

$bowlingCodes = @{
    2  = "4235"; 3  = "4533"; 4  = "4535"; 5  = "4536"; 6  = "4618";
    7  = "4621"; 8  = "4623"; 9  = "4624"; 10 = "4637"; 11 = "4640";
    12 = "4656"; 13 = "4657"; 14 = "4658"; 15 = "4679"; 16 = "4682";
    17 = "4685"; 18 = "4687"; 19 = "4689"; 20 = "4691"; 21 = "4692";
    22 = "4695"; 23 = "4725"; 24 = "4728"; 25 = "4732"
}
$bowlingSheet.Range("B2:B25").NumberFormat = "@"
foreach ($r in $bowlingCodes.Keys) {
    $bowlingSheet.Cells.Item($r, 2).Value = $bowlingCodes[$r]
}

# ----------------------------------------------------------------------------
# 5. Populate "ODI Batting Extra"
# ----------------------------------------------------------------------------
$extraSheet = $wb.Worksheets.Item("ODI Batting Extra")

$exHeader = $extraSheet.Range("A1:F1")
Format-HeaderRange $exHeader
$extraSheet.Cells.Item(1,1).Value = "MATCH_CODE"
$extraSheet.Cells.Item(1,2).Value = "BATTING_POSITION"
$extraSheet.Cells.Item(1,3).Value = "NUM_4"
$extraSheet.Cells.Item(1,4).Value = "NUM_6"
$extraSheet.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$extraSheet.Cells.Item(1,6).Value = "MAN_OF_MATCH"

# Column A (MATCH_CODE) holds its numeric-looking values as text, same as
# the rest of the workbook's scraped columns; columns C/D/E get the same
# text treatment but only on the cells that actually receive a value below
# (so untouched rows don't grow spurious blank-but-styled cells).
$extraSheet.Range("A2:A21").NumberFormat = "@"

# row -> (MATCH_CODE, BATTING_POSITION[number-or-blank], NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH)
$extraRows = @(
    @(2,  "4618", 9,    "0", "0", "0.32%", "NO"),
    @(3,  "4621", $null, $null, $null, $null, "NO"),
    @(4,  "4623", 10,   $null, $null, $null, "NO"),
    @(5,  "4624", $null, $null, $null, $null, "NO"),
    @(6,  "4637", 10,   "1", "0", "4.84%", "NO"),
    @(7,  "4640", 10,   "0", "0", "0.75%", "NO"),
    @(8,  "4656", $null, $null, $null, $null, "NO"),
    @(9,  "4657", $null, $null, $null, $null, "NO"),
    @(10, "4658", 10,   $null, $null, $null, "NO"),
    @(11, "4679", $null, $null, $null, $null, "NO"),
    @(12, "4682", $null, $null, $null, $null, "NO"),
    @(13, "4685", 10,   $null, $null, $null, "NO"),
    @(14, "4687", 11,   "0", "0", $null, "NO"),
    @(15, "4689", $null, $null, $null, $null, $null),
    @(16, "4691", $null, $null, $null, $null, $null),
    @(17, "4692", $null, $null, $null, $null, $null),
    @(18, "4695", $null, $null, $null, $null, $null),
    @(19, "4725", $null, $null, $null, $null, $null),
    @(20, "4728", $null, $null, $null, $null, $null),
    @(21, "4732", $null, $null, $null, $null, $null)
)

foreach ($entry in $extraRows) {
    $r = $entry[0]
    $matchCode = $entry[1]
    $battingPos = $entry[2]
    $num4 = $entry[3]
    $num6 = $entry[4]
    $pct = $entry[5]
    $mom = $entry[6]

    $extraSheet.Cells.Item($r, 1).Value = $matchCode
    if ($null -ne $battingPos) {
        $extraSheet.Cells.Item($r, 2).Value = $battingPos
    }
    if ($null -ne $num4) {
        $extraSheet.Cells.Item($r, 3).NumberFormat = "@"
        $extraSheet.Cells.Item($r, 3).Value = $num4
    }
    if ($null -ne $num6) {
        $extraSheet.Cells.Item($r, 4).NumberFormat = "@"
        $extraSheet.Cells.Item($r, 4).Value = $num6
    }
    if ($null -ne $pct) {
        $extraSheet.Cells.Item($r, 5).NumberFormat = "@"
        $extraSheet.Cells.Item($r, 5).Value = $pct
    }
    if ($null -ne $mom) {
        $extraSheet.Cells.Item($r, 6).Value = $mom
    }
}

# ----------------------------------------------------------------------------
# Leave the workbook focused on the first sheet, matching activeTab="0".
# ----------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Activate()
